$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("A8").Value = "Profile"
$ws.Range("B8").Value = "Profile Info Page"
$ws.Range("C8").Value = "YES"

$ws.Range("C12").Select()
